$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 2579.7837  # H70: 4044.2354 -> 2579.7837
$ws.Cells.Item(70, 9).Value = 7993.143  # I70: 7206.5 -> 7993.143
$ws.Cells.Item(70, 10).Value = 1316.6666  # J70: 1233.3334 -> 1316.6666
$ws.Cells.Item(70, 11).Value = 23979.429  # K70: 21619.5 -> 23979.429
$ws.Cells.Item(70, 12).Value = 3949.9998  # L70: 3700.0002 -> 3949.9998
$ws.Cells.Item(70, 13).Value = -23709.429  # M70: -21349.5 -> -23709.429
$ws.Cells.Item(70, 14).Value = -4489.9998  # N70: -4240.0002 -> -4489.9998
$ws.Cells.Item(73, 8).Value = 2579.7837  # H73: 4044.2354 -> 2579.7837
$ws.Cells.Item(73, 9).Value = 7993.143  # I73: 7206.5 -> 7993.143
$ws.Cells.Item(73, 10).Value = 1316.6666  # J73: 1233.3334 -> 1316.6666
$ws.Cells.Item(73, 11).Value = 23979.429  # K73: 21619.5 -> 23979.429
$ws.Cells.Item(73, 12).Value = 3949.9998  # L73: 3700.0002 -> 3949.9998
$ws.Cells.Item(73, 13).Value = -23043.429  # M73: -20683.5 -> -23043.429
$ws.Cells.Item(73, 14).Value = -5821.9998  # N73: -5572.0002 -> -5821.9998
$ws.Cells.Item(74, 8).Value = 4309.1714  # H74: 5374.6665 -> 4309.1714
$ws.Cells.Item(74, 9).Value = 4455.3076  # I74: 5205 -> 4455.3076
$ws.Cells.Item(74, 10).Value = 4222.8184  # J74: 5568.5713 -> 4222.8184
$ws.Cells.Item(74, 11).Value = 4455.3076  # K74: 5205 -> 4455.3076
$ws.Cells.Item(74, 12).Value = 4222.8184  # L74: 5568.5713 -> 4222.8184
$ws.Cells.Item(74, 13).Value = -3519.3076  # M74: -4269 -> -3519.3076
$ws.Cells.Item(74, 14).Value = -6094.8184  # N74: -7440.5713 -> -6094.8184
$ws.Cells.Item(75, 8).Value = 22611.4  # H75: 26462.8 -> 22611.4
$ws.Cells.Item(75, 10).Value = 22611.4  # J75: 26462.8 -> 22611.4
$ws.Cells.Item(75, 12).Value = 22611.4  # L75: 26462.8 -> 22611.4
$ws.Cells.Item(75, 14).Value = -24483.4  # N75: -28334.8 -> -24483.4
$ws.Cells.Item(77, 8).Value = 4309.1714  # H77: 5374.6665 -> 4309.1714
$ws.Cells.Item(77, 9).Value = 4455.3076  # I77: 5205 -> 4455.3076
$ws.Cells.Item(77, 10).Value = 4222.8184  # J77: 5568.5713 -> 4222.8184
$ws.Cells.Item(77, 11).Value = 22276.538  # K77: 26025 -> 22276.538
$ws.Cells.Item(77, 12).Value = 21114.092  # L77: 27842.8565 -> 21114.092
$ws.Cells.Item(77, 13).Value = -17596.538  # M77: -21345 -> -17596.538
$ws.Cells.Item(77, 14).Value = -30474.092  # N77: -37202.85649999999 -> -30474.092
$ws.Cells.Item(78, 8).Value = 22611.4  # H78: 26462.8 -> 22611.4
$ws.Cells.Item(78, 10).Value = 22611.4  # J78: 26462.8 -> 22611.4
$ws.Cells.Item(78, 12).Value = 67834.20000000001  # L78: 79388.39999999999 -> 67834.20000000001
$ws.Cells.Item(78, 14).Value = -77194.20000000001  # N78: -88748.39999999999 -> -77194.20000000001
$ws.Cells.Item(81, 8).Value = 0  # H81: 6000 -> 0
$ws.Cells.Item(81, 9).Value = 0  # I81: 6000 -> 0
$ws.Cells.Item(81, 11).Value = 0  # K81: 6000 -> 0
$ws.Cells.Item(81, 13).Value = $null  # M81: -5002 -> (removed)
$ws.Cells.Item(84, 8).Value = 0  # H84: 6000 -> 0
$ws.Cells.Item(84, 9).Value = 0  # I84: 6000 -> 0
$ws.Cells.Item(84, 11).Value = 0  # K84: 18000 -> 0
$ws.Cells.Item(84, 13).Value = $null  # M84: -13008 -> (removed)

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(60, 8).Value = 0  # H60: 27800 -> 0
$ws.Cells.Item(60, 10).Value = 0  # J60: 27800 -> 0
$ws.Cells.Item(60, 12).Value = 0  # L60: 27800 -> 0
$ws.Cells.Item(60, 14).Value = $null  # N60: -29266 -> (removed)
$ws.Cells.Item(97, 8).Value = 1177.5186  # H97: 972.7353000000001 -> 1177.5186
$ws.Cells.Item(97, 9).Value = 1184.3462  # I97: 971.9091 -> 1184.3462
$ws.Cells.Item(97, 11).Value = 1184.3462  # K97: 971.9091 -> 1184.3462
$ws.Cells.Item(97, 13).Value = -688.3462  # M97: -475.9091 -> -688.3462

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1264.3077  # H94: 1149.0667 -> 1264.3077
$ws.Cells.Item(94, 9).Value = 674.2  # I94: 646.8570999999999 -> 674.2
$ws.Cells.Item(94, 10).Value = 3231.3333  # J94: 2320.889 -> 3231.3333
$ws.Cells.Item(94, 11).Value = 674.2  # K94: 646.8570999999999 -> 674.2
$ws.Cells.Item(94, 12).Value = 3231.3333  # L94: 2320.889 -> 3231.3333
$ws.Cells.Item(94, 13).Value = -223.2  # M94: -195.8570999999999 -> -223.2
$ws.Cells.Item(94, 14).Value = -4133.3333  # N94: -3222.889 -> -4133.3333

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(87, 8).Value = 133553.33  # H87: 50000 -> 133553.33
$ws.Cells.Item(87, 10).Value = 133553.33  # J87: 50000 -> 133553.33
$ws.Cells.Item(87, 12).Value = 133553.33  # L87: 50000 -> 133553.33
$ws.Cells.Item(87, 14).Value = -135925.33  # N87: -52372 -> -135925.33
$ws.Cells.Item(90, 8).Value = 133553.33  # H90: 50000 -> 133553.33
$ws.Cells.Item(90, 10).Value = 133553.33  # J90: 50000 -> 133553.33
$ws.Cells.Item(90, 12).Value = 400659.99  # L90: 150000 -> 400659.99
$ws.Cells.Item(90, 14).Value = -412515.99  # N90: -161856 -> -412515.99
$ws.Cells.Item(134, 8).Value = 2172.875  # H134: 2213.9355 -> 2172.875
$ws.Cells.Item(134, 9).Value = 1451.5  # I134: 1488.2667 -> 1451.5
$ws.Cells.Item(134, 11).Value = 4354.5  # K134: 4464.800099999999 -> 4354.5
$ws.Cells.Item(134, 13).Value = -1819.5  # M134: -1929.800099999999 -> -1819.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 702.4286  # H5: 872.95 -> 702.4286
$ws.Cells.Item(5, 9).Value = 641.0769  # I5: 644.875 -> 641.0769
$ws.Cells.Item(5, 10).Value = 1500  # J5: 1025 -> 1500
$ws.Cells.Item(5, 11).Value = 1923.2307  # K5: 1934.625 -> 1923.2307
$ws.Cells.Item(5, 12).Value = 4500  # L5: 3075 -> 4500
$ws.Cells.Item(5, 13).Value = -1811.2307  # M5: -1822.625 -> -1811.2307
$ws.Cells.Item(5, 14).Value = -4724  # N5: -3299 -> -4724
$ws.Cells.Item(82, 8).Value = 6142.857  # H82: 5237.5 -> 6142.857
$ws.Cells.Item(82, 9).Value = 2000  # I82: 950 -> 2000
$ws.Cells.Item(82, 10).Value = 6833.3335  # J82: 6666.6665 -> 6833.3335
$ws.Cells.Item(82, 11).Value = 6000  # K82: 2850 -> 6000
$ws.Cells.Item(82, 12).Value = 20500.0005  # L82: 19999.9995 -> 20500.0005
$ws.Cells.Item(82, 13).Value = -5594  # M82: -2444 -> -5594
$ws.Cells.Item(82, 14).Value = -21312.0005  # N82: -20811.9995 -> -21312.0005
$ws.Cells.Item(85, 8).Value = 6142.857  # H85: 5237.5 -> 6142.857
$ws.Cells.Item(85, 9).Value = 2000  # I85: 950 -> 2000
$ws.Cells.Item(85, 10).Value = 6833.3335  # J85: 6666.6665 -> 6833.3335
$ws.Cells.Item(85, 11).Value = 6000  # K85: 2850 -> 6000
$ws.Cells.Item(85, 12).Value = 20500.0005  # L85: 19999.9995 -> 20500.0005
$ws.Cells.Item(85, 13).Value = -4596  # M85: -1446 -> -4596
$ws.Cells.Item(85, 14).Value = -23308.0005  # N85: -22807.9995 -> -23308.0005
$ws.Cells.Item(98, 8).Value = 10760  # H98: 9131.333000000001 -> 10760
$ws.Cells.Item(98, 9).Value = 933.3333  # I98: 947 -> 933.3333
$ws.Cells.Item(98, 11).Value = 2799.9999  # K98: 2841 -> 2799.9999
$ws.Cells.Item(98, 13).Value = -1301.9999  # M98: -1343 -> -1301.9999
$ws.Cells.Item(135, 8).Value = 702.4286  # H135: 872.95 -> 702.4286
$ws.Cells.Item(135, 9).Value = 641.0769  # I135: 644.875 -> 641.0769
$ws.Cells.Item(135, 10).Value = 1500  # J135: 1025 -> 1500
$ws.Cells.Item(135, 11).Value = 5769.6921  # K135: 5803.875 -> 5769.6921
$ws.Cells.Item(135, 12).Value = 13500  # L135: 9225 -> 13500
$ws.Cells.Item(135, 13).Value = -3234.6921  # M135: -3268.875 -> -3234.6921
$ws.Cells.Item(135, 14).Value = -18570  # N135: -14295 -> -18570

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(53, 8).Value = 3965.8333  # H53: 0 -> 3965.8333
$ws.Cells.Item(53, 9).Value = 759  # I53: 0 -> 759
$ws.Cells.Item(53, 10).Value = 20000  # J53: 0 -> 20000
$ws.Cells.Item(53, 11).Value = 759  # K53: 0 -> 759
$ws.Cells.Item(53, 12).Value = 20000  # L53: 0 -> 20000
$ws.Cells.Item(53, 13).Value = -128  # M53: None -> -128
$ws.Cells.Item(53, 14).Value = -21262  # N53: None -> -21262
$ws.Cells.Item(80, 8).Value = 6355.5557  # H80: 5483.3335 -> 6355.5557
$ws.Cells.Item(80, 9).Value = 8191.6665  # I80: 7761.5386 -> 8191.6665
$ws.Cells.Item(80, 10).Value = 2683.3333  # J80: 2790.9092 -> 2683.3333
$ws.Cells.Item(80, 11).Value = 8191.6665  # K80: 7761.5386 -> 8191.6665
$ws.Cells.Item(80, 12).Value = 2683.3333  # L80: 2790.9092 -> 2683.3333
$ws.Cells.Item(80, 13).Value = -7193.6665  # M80: -6763.5386 -> -7193.6665
$ws.Cells.Item(80, 14).Value = -4679.3333  # N80: -4786.9092 -> -4679.3333
$ws.Cells.Item(83, 8).Value = 6355.5557  # H83: 5483.3335 -> 6355.5557
$ws.Cells.Item(83, 9).Value = 8191.6665  # I83: 7761.5386 -> 8191.6665
$ws.Cells.Item(83, 10).Value = 2683.3333  # J83: 2790.9092 -> 2683.3333
$ws.Cells.Item(83, 11).Value = 40958.3325  # K83: 38807.693 -> 40958.3325
$ws.Cells.Item(83, 12).Value = 13416.6665  # L83: 13954.546 -> 13416.6665
$ws.Cells.Item(83, 13).Value = -35966.3325  # M83: -33815.693 -> -35966.3325
$ws.Cells.Item(83, 14).Value = -23400.6665  # N83: -23938.546 -> -23400.6665
$ws.Cells.Item(97, 8).Value = 936.41174  # H97: 1017.53845 -> 936.41174
$ws.Cells.Item(97, 9).Value = 860.5333000000001  # I97: 935.6667 -> 860.5333000000001
$ws.Cells.Item(97, 10).Value = 1505.5  # J97: 2000 -> 1505.5
$ws.Cells.Item(97, 11).Value = 860.5333000000001  # K97: 935.6667 -> 860.5333000000001
$ws.Cells.Item(97, 12).Value = 1505.5  # L97: 2000 -> 1505.5
$ws.Cells.Item(97, 13).Value = -364.5333000000001  # M97: -439.6667 -> -364.5333000000001
$ws.Cells.Item(97, 14).Value = -2497.5  # N97: -2992 -> -2497.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(74, 8).Value = 28205  # H74: 33665 -> 28205
$ws.Cells.Item(74, 10).Value = 28205  # J74: 33665 -> 28205
$ws.Cells.Item(74, 12).Value = 28205  # L74: 33665 -> 28205
$ws.Cells.Item(74, 14).Value = -30201  # N74: -35661 -> -30201
$ws.Cells.Item(77, 8).Value = 28205  # H77: 33665 -> 28205
$ws.Cells.Item(77, 10).Value = 28205  # J77: 33665 -> 28205
$ws.Cells.Item(77, 12).Value = 84615  # L77: 100995 -> 84615
$ws.Cells.Item(77, 14).Value = -94599  # N77: -110979 -> -94599
$ws.Cells.Item(82, 8).Value = 1877.35  # H82: 1985.9445 -> 1877.35
$ws.Cells.Item(82, 9).Value = 1800  # I82: 2700 -> 1800
$ws.Cells.Item(82, 11).Value = 1800  # K82: 2700 -> 1800
$ws.Cells.Item(82, 13).Value = -1439  # M82: -2339 -> -1439
$ws.Cells.Item(85, 8).Value = 1877.35  # H85: 1985.9445 -> 1877.35
$ws.Cells.Item(85, 9).Value = 1800  # I85: 2700 -> 1800
$ws.Cells.Item(85, 11).Value = 1800  # K85: 2700 -> 1800
$ws.Cells.Item(85, 13).Value = -552  # M85: -1452 -> -552

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(57, 8).Value = 37000  # H57: 0 -> 37000
$ws.Cells.Item(57, 10).Value = 37000  # J57: 0 -> 37000
$ws.Cells.Item(57, 12).Value = 37000  # L57: 0 -> 37000
$ws.Cells.Item(57, 14).Value = -38508  # N57: None -> -38508
$ws.Cells.Item(81, 8).Value = 2017.3077  # H81: 1929.1666 -> 2017.3077
$ws.Cells.Item(81, 10).Value = 3791.6667  # J81: 4150 -> 3791.6667
$ws.Cells.Item(81, 12).Value = 7583.3334  # L81: 8300 -> 7583.3334
$ws.Cells.Item(81, 14).Value = -9705.3334  # N81: -10422 -> -9705.3334
$ws.Cells.Item(84, 8).Value = 2017.3077  # H84: 1929.1666 -> 2017.3077
$ws.Cells.Item(84, 10).Value = 3791.6667  # J84: 4150 -> 3791.6667
$ws.Cells.Item(84, 12).Value = 37916.667  # L84: 41500 -> 37916.667
$ws.Cells.Item(84, 14).Value = -48524.667  # N84: -52108 -> -48524.667
$ws.Cells.Item(136, 8).Value = 22240746  # H136: 23411364 -> 22240746
$ws.Cells.Item(136, 9).Value = 25667960  # I136: 29442636 -> 25667960
$ws.Cells.Item(136, 10).Value = 15875921  # J136: 14495571 -> 15875921
$ws.Cells.Item(136, 11).Value = 77003880  # K136: 88327908 -> 77003880
$ws.Cells.Item(136, 12).Value = 47627763  # L136: 43486713 -> 47627763
$ws.Cells.Item(136, 13).Value = -77001330  # M136: -88325358 -> -77001330
$ws.Cells.Item(136, 14).Value = -47632863  # N136: -43491813 -> -47632863
